# Update "想去人数" (interested-count) figures that changed between scrapes.
# Same six events are listed on both the "展览" (Exhibition) sheet and the
# "全部类型" (All types) roll-up sheet, so each value bumps in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value, for the "展览" sheet (column F)
$exhibitUpdates = @{
    4  = 3376
    6  = 4871
    20 = 4778
    23 = 10
    27 = 249
    36 = 801
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new value, for the "全部类型" sheet (column F)
$allUpdates = @{
    8  = 3376
    10 = 4871
    25 = 4778
    28 = 10
    32 = 249
    42 = 801
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
